$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.481.08"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "1.674.09"
$ws.Range("E3").Value = "  +1.65%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5323"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.17%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +3.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06387"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.74"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07799"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.674.86"
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.493"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5569"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "0.0₅8337"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.96%  "
$ws.Range("D17").Value = "26.503.65"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.758"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.312"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1277"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.410"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.426"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06256"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.64%  "
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("E31").Value = "  +5.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.426"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.692"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.010"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6147"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.421"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.783"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.135"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01617"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "1.093.00"
$ws.Range("E40").Value = "  +6.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8605"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").Value = "1.820.48"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "58.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.180"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.519"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05196"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.008"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.15%  "
